# SAV-700: Update charts importer test fixture
#  - Ensure ID of charting date recorded program data element is constant
#    ("Test Chart" sheet: testchartcode0 -> PatientChartingDate)
#  - Ensure ID of all complex chart core questions is constant
#    ("Core" sheet: testchartcorecodeN -> same value as the question's `type`
#    column, i.e. the code now mirrors column B)

$wb = $excel.ActiveWorkbook

$core = $wb.Worksheets.Item("Core")
$chart = $wb.Worksheets.Item("Test Chart")

# Grab the (slightly) distinct cell format that currently lives on
# Core!R5 ("historical") before it gets normalised below, so it can be
# carried over to the cell that now needs to stand out instead
# (Test Chart!A2, the constant charting-date code).
$core.Range("R5").Copy()
$chart.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

# --- Core sheet: make the `code` column (A) match the `type` column (B) ---
$core.Range("A2").Value = $core.Range("B2").Value2
$core.Range("A3").Value = $core.Range("B3").Value2
$core.Range("A4").Value = $core.Range("B4").Value2
$core.Range("A5").Value = $core.Range("B5").Value2

# R5's cell format is normalised back to match the rest of the column
# (R2:R4) now that it no longer needs to stand out.
$core.Range("R2").Copy()
$core.Range("R5").PasteSpecial(-4122)  # xlPasteFormats

# --- Test Chart sheet: rename the charting-date question's constant code ---
$chart.Range("A2").Value = "PatientChartingDate"

$excel.CutCopyMode = 0
